# Applies the "Started progress on final Clemens comments" edit:
#  - updates the contrast data table (V2:X10) on Sheet1
#  - repositions/resizes the embedded chart
#  - updates the active selection

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# -----------------------------------------------------------------
# 1. Update the data values that feed the chart (columns V, W, X,
#    rows 2-10). Literals below are the shortest decimal strings that
#    round-trip to the target doubles.
# -----------------------------------------------------------------
$ws.Range("V2").Value = 0.6
$ws.Range("W2").Value = 0.78
$ws.Range("X2").Value = 0.56

$ws.Range("V3").Value = 0.6
$ws.Range("W3").Value = 0.6
$ws.Range("X3").Value = 0.64

$ws.Range("V4").Value = 0.13
$ws.Range("W4").Value = 0.2
$ws.Range("X4").Value = 0.14

$ws.Range("V5").Value = 0.8438
$ws.Range("W5").Value = 0.9968
$ws.Range("X5").Value = 0.9406

$ws.Range("V6").Value = 0.8946
$ws.Range("W6").Value = 0.9757
$ws.Range("X6").Value = 0.968

$ws.Range("V7").Value = 0.8027
$ws.Range("W7").Value = 1
$ws.Range("X7").Value = 0.9528

$ws.Range("V8").Value = 0.78
$ws.Range("W8").Value = 0.9968
$ws.Range("X8").Value = 0.94

$ws.Range("V9").Value = 0.67
$ws.Range("W9").Value = 0.9757
$ws.Range("X9").Value = 0.81

$ws.Range("V10").Value = 0.79
$ws.Range("W10").Value = 1
$ws.Range("X10").Value = 0.95

# -----------------------------------------------------------------
# 2. Move/resize the chart so its anchor matches:
#      from col=3 (0-based) offset 333375 EMU, row=3 offset 61911 EMU
#      to   col=18           offset 85725  EMU, row=30 offset 104774 EMU
#    Convert the target column/row + offset position into points
#    (Excel's Left/Top/Width/Height units) using the sheet's actual
#    column widths / row heights, since no custom widths are set.
# -----------------------------------------------------------------
$EMU_PER_POINT = 12700

$leftAccum = 0
for ($i = 1; $i -le 3; $i++) {
    $leftAccum = $leftAccum + $ws.Columns.Item($i).Width
}
$newLeft = $leftAccum + (333375 / $EMU_PER_POINT)

$topAccum = 0
for ($i = 1; $i -le 3; $i++) {
    $topAccum = $topAccum + $ws.Rows.Item($i).Height
}
$newTop = $topAccum + (61911 / $EMU_PER_POINT)

$rightAccum = 0
for ($i = 1; $i -le 18; $i++) {
    $rightAccum = $rightAccum + $ws.Columns.Item($i).Width
}
$newRight = $rightAccum + (85725 / $EMU_PER_POINT)

$bottomAccum = 0
for ($i = 1; $i -le 30; $i++) {
    $bottomAccum = $bottomAccum + $ws.Rows.Item($i).Height
}
$newBottom = $bottomAccum + (104774 / $EMU_PER_POINT)

$newWidth = $newRight - $newLeft
$newHeight = $newBottom - $newTop

$co = $ws.ChartObjects().Item(1)
$co.Left = $newLeft
$co.Top = $newTop
$co.Width = $newWidth
$co.Height = $newHeight

# -----------------------------------------------------------------
# 3. Update the active selection shown when the sheet is reopened.
# -----------------------------------------------------------------
$ws.Range("T13").Select()
